$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Concepts": insert BAM/BAI rows after CRAI, append FASTQ row at end.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Concepts")

# Shift VCF..BED (rows 4-9) down by two rows (-> rows 6-11), keeping format.
$ws.Rows("4:5").Insert(-4121)

# Clone row 11 (BED | BED File) down into new row 12 (values then formats,
# so the copy keeps the shared-string text type and the "s=2" style).
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial(-4163)
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)

# The two rows that Insert() opened up (4 and 5) are blank with a fresh
# style. Clone row 6 (still intact: Level=1, style 2) into both so they
# pick up the "1" level value, the correct style, and an empty D cell.
$ws.Range("A6:D6").Copy()
$ws.Range("A4:D4").PasteSpecial(-4163)
$ws.Range("A6:D6").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

$ws.Range("A6:D6").Copy()
$ws.Range("A5:D5").PasteSpecial(-4163)
$ws.Range("A6:D6").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# Fill in the actual new concept rows.
$ws.Range("B4").Value = "BAM"
$ws.Range("C4").Value = "BAM File"
$ws.Range("B5").Value = "BAI"
$ws.Range("C5").Value = "BAI File"
$ws.Range("B12").Value = "FASTQ"
$ws.Range("C12").Value = "FASTQ File"

# ---------------------------------------------------------------------------
# Sheet "Metadata": refresh Date and Count.
# ---------------------------------------------------------------------------
$ms = $wb.Worksheets.Item("Metadata")

$ms.Range("B8").Value = "2025-05-21T20:08:08+00:00"

# "11" would otherwise be auto-typed as a number; force it back to text
# (matching the original shared-string cell) and then re-apply the
# surrounding normal-row style so it collapses back onto the shared "s=2"
# style instead of forking a new number-format style.
$ms.Range("B22").NumberFormat = "@"
$ms.Range("B22").Value = "11"
$ms.Range("B21").Copy()
$ms.Range("B22").PasteSpecial(-4122)
